# Added Texture coords and index
# Populate the newly-added texture-coordinate columns (AO/AP) for rows 4-35
# and the single index value in AT21, matching the "Added Texture coords
# and index" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO4").Value = 0.5625
$ws.Range("AP4").Value = 0.5

$ws.Range("AO5").Value = 0.53125
$ws.Range("AP5").Value = 0.53125

$ws.Range("AO6").Value = 0.5302734375
$ws.Range("AP6").Value = 0.5615234375

$ws.Range("AO7").Value = 0.5625
$ws.Range("AP7").Value = 0.59375

$ws.Range("AO8").Value = 0.53125
$ws.Range("AP8").Value = 0.625

$ws.Range("AO9").Value = 0.5
$ws.Range("AP9").Value = 0.59375

$ws.Range("AO10").Value = 0.5
$ws.Range("AP10").Value = 0.5615234375

$ws.Range("AO11").Value = 0.5
$ws.Range("AP11").Value = 0.53125

$ws.Range("AO12").Value = 0.5
$ws.Range("AP12").Value = 0.5

$ws.Range("AO13").Value = 0.5302734375
$ws.Range("AP13").Value = 0.46875

$ws.Range("AO14").Value = 0.6552734375
$ws.Range("AP14").Value = 0.46875

$ws.Range("AO15").Value = 0.6552734375
$ws.Range("AP15").Value = 0.40625

$ws.Range("AO16").Value = 0.6875
$ws.Range("AP16").Value = 0.40625

$ws.Range("AO17").Value = 0.6875
$ws.Range("AP17").Value = 0.5615234375

$ws.Range("AO18").Value = 0.6552734375
$ws.Range("AP18").Value = 0.5615234375

$ws.Range("AO19").Value = 0.6552734375
$ws.Range("AP19").Value = 0.5

$ws.Range("AO20").Value = 0.375
$ws.Range("AP20").Value = 0.5

$ws.Range("AO21").Value = 0.40625
$ws.Range("AP21").Value = 0.53125
$ws.Range("AT21").Value = 1024

$ws.Range("AO22").Value = 0.40625
$ws.Range("AP22").Value = 0.5615234375

$ws.Range("AO23").Value = 0.375
$ws.Range("AP23").Value = 0.59375

$ws.Range("AO24").Value = 0.40625
$ws.Range("AP24").Value = 0.625

$ws.Range("AO25").Value = 0.4375
$ws.Range("AP25").Value = 0.59375

$ws.Range("AO26").Value = 0.4375
$ws.Range("AP26").Value = 0.5615234375

$ws.Range("AO27").Value = 0.4375
$ws.Range("AP27").Value = 0.53125

$ws.Range("AO28").Value = 0.22265625
$ws.Range("AP28").Value = 0.5

$ws.Range("AO29").Value = 0.40625
$ws.Range("AP29").Value = 0.46875

$ws.Range("AO30").Value = 0.28125
$ws.Range("AP30").Value = 0.46875

$ws.Range("AO31").Value = 0.28125
$ws.Range("AP31").Value = 0.40625

$ws.Range("AO32").Value = 0.2490234375
$ws.Range("AP32").Value = 0.40625

$ws.Range("AO33").Value = 0.2490234375
$ws.Range("AP33").Value = 0.5615234375

$ws.Range("AO34").Value = 0.28125
$ws.Range("AP34").Value = 0.5615234375

$ws.Range("AO35").Value = 0.28125
$ws.Range("AP35").Value = 0.5

# Mirror the author's final selection/view state (scroll so O7 is the
# top-left cell, then select U3:U34 as the last user action).
$excel.Goto($ws.Range("O7"), $true)
$ws.Range("U3:U34").Select()
